# "Generate Report for Handback"
#
# The handback xliffs for the 3961d949-...md source file have now come
# back in sync with en-US for both target locales (zh-cn, de-de). This
# script fills in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on the two locale detail sheets for
# that row, links the target file name back to the source doc on GitHub
# (same link used by column A), and flips the row's Status text (which is
# shared with the Overview sheet's per-locale status cells) from
# "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$sourceMdName = "3961d949-0846-4b9d-904d-bbdf85b459ed.md"
$sourceMdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/db52bba40c3c77ecd9cac36b142ec1a5919ae73b/e2e/3961d949-0846-4b9d-904d-bbdf85b459ed.md"
$statusText   = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: the zh-cn / de-de status cells for this row pick up the
# same new status text.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

# Columns widened to fit the longer status text.
$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Range("I2").Value = $sourceMdName
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $sourceMdUrl, "", "", $sourceMdName) | Out-Null
$zhcn.Range("J2").Value = "3961d949-0846-4b9d-904d-bbdf85b459ed.4b8f675dd86eba16a9241862801644053b614db6.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-10-20 00:43:01"

$zhcn.Range("I3").Value = $sourceMdName
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $sourceMdUrl, "", "", $sourceMdName) | Out-Null
$zhcn.Range("J3").Value = "3961d949-0846-4b9d-904d-bbdf85b459ed.4b8f675dd86eba16a9241862801644053b614db6.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-10-20 00:43:01"

# Columns widened to fit the longer status text / new file names.
$zhcn.Columns.Item(3).ColumnWidth = 29.15
$zhcn.Columns.Item(9).ColumnWidth = 39.15
$zhcn.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Range("I2").Value = $sourceMdName
$dede.Hyperlinks.Add($dede.Range("I2"), $sourceMdUrl, "", "", $sourceMdName) | Out-Null
$dede.Range("J2").Value = "3961d949-0846-4b9d-904d-bbdf85b459ed.4b8f675dd86eba16a9241862801644053b614db6.de-de.xlf"
$dede.Range("K2").Value = "2016-10-20 00:43:20"

$dede.Range("I3").Value = $sourceMdName
$dede.Hyperlinks.Add($dede.Range("I3"), $sourceMdUrl, "", "", $sourceMdName) | Out-Null
$dede.Range("J3").Value = "3961d949-0846-4b9d-904d-bbdf85b459ed.4b8f675dd86eba16a9241862801644053b614db6.de-de.xlf"
$dede.Range("K3").Value = "2016-10-20 00:43:20"

# Columns widened to fit the longer status text / new file names.
$dede.Columns.Item(3).ColumnWidth = 29.15
$dede.Columns.Item(9).ColumnWidth = 39.15
$dede.Columns.Item(10).ColumnWidth = 39.15

Write-Host "Handback report generated."
